# Updated cryptos list on Fri May 12 10:54:55 UTC 2023 with GitHub Actions
#
# Refresh the live "Price" (D) / "Volume(1h)" (E) columns, and re-sort a few
# rows (46-49) whose ranking order changed between scrapes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source feed writes plain numeric-looking strings (e.g. "1.002") into the
# Price column as literal text, not numbers. Force Text format first so Excel
# does not re-interpret them as numbers on assignment (multi-dot values such as
# "26.447.63" are left alone below since Excel cannot parse them as numbers anyway).

$ws.Range("D2").Value = "26.447.63"
$ws.Range("E2").Value = "  -3.84%  "

$ws.Range("D3").Value = "1.773.56"
$ws.Range("E3").Value = "  -2.89%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.76"
$ws.Range("E6").Value = "  -1.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4308"
$ws.Range("E7").Value = "  +1.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3671"
$ws.Range("E8").Value = "  +1.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07238"
$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8506"
$ws.Range("E10").Value = "  -1.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.33"
$ws.Range("E11").Value = "  -1.50%  "

$ws.Range("D12").Value = "1.780.09"
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.448"
$ws.Range("E13").Value = "  -0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.250"
$ws.Range("E14").Value = "  -2.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06861"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.64"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008696"
$ws.Range("E18").Value = "  -2.57%  "

$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("E20").Value = "  -1.90%  "

$ws.Range("D21").Value = "26.448.08"
$ws.Range("E21").Value = "  -3.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.118"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.31"
$ws.Range("E23").Value = "  +3.87%  "

$ws.Range("D24").Value = "1.988.94"
$ws.Range("E24").Value = "  -2.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.35"
$ws.Range("E25").Value = "  -1.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.857"
$ws.Range("E26").Value = "  -6.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.20"
$ws.Range("E27").Value = "  -2.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.094"
$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.77"
$ws.Range("E29").Value = "  +0.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.720"
$ws.Range("E30").Value = "  -4.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08962"
$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7256"
$ws.Range("E32").Value = "  -3.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.119"
$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.334"
$ws.Range("E34").Value = "  -4.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.750"
$ws.Range("E35").Value = "  -7.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.080"
$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05169"
$ws.Range("E38").Value = "  -2.21%  "

$ws.Range("E39").Value = "  -1.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4940"
$ws.Range("E40").Value = "  -2.60%  "

$ws.Range("E41").Value = "  -2.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.534"
$ws.Range("E42").Value = "  -9.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.229"
$ws.Range("E43").Value = "  -3.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.047"
$ws.Range("E44").Value = "  -3.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.00"
$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.584"
$ws.Range("E50").Value = "  -2.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.746"
$ws.Range("E51").Value = "  +2.00%  "

# --- Rows 46-49 reordered (ranking shuffle) with refreshed coin/link/price/volume ---
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.18"
$ws.Range("E46").Value = "  -2.92%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4508"
$ws.Range("E48").Value = "  -3.57%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06198"
$ws.Range("E49").Value = "  -4.13%  "

